$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.245.74'
$ws.Range('E2').Value = '  -1.49%  '
$ws.Range('D3').Value = '3.520.02'
$ws.Range('E3').Value = '  -0.22%  '
$ws.Range('E4').Value = '  +0.26%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '585.32'
$ws.Range('E5').Value = '  +5.83%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '178.81'
$ws.Range('E6').Value = '  -6.60%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.630'
$ws.Range('E7').Value = '  +3.58%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  +0.83%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.163'
$ws.Range('E10').Value = '  +6.97%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '56.08'
$ws.Range('E11').Value = '  +1.65%  '
$ws.Range('E12').Value = '  +3.74%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.30'
$ws.Range('E13').Value = '  -0.76%  '
$ws.Range('D14').Value = '4.084.05'
$ws.Range('E14').Value = '  +0.36%  '
$ws.Range('D15').Value = '3.534.10'
$ws.Range('E15').Value = '  +0.52%  '
$ws.Range('E16').Value = '  +0.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.41'
$ws.Range('E17').Value = '  +0.99%  '
$ws.Range('D18').Value = '66.250.52'
$ws.Range('E18').Value = '  -1.39%  '
$ws.Range('E19').Value = '  +1.41%  '
$ws.Range('E20').Value = '  +1.37%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '416.61'
$ws.Range('E21').Value = '  -2.71%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.33'
$ws.Range('E22').Value = '  +10.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.42'
$ws.Range('E23').Value = '  +5.54%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.73'
$ws.Range('E24').Value = '  +0.98%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.39'
$ws.Range('E25').Value = '  +11.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.10'
$ws.Range('E26').Value = '  -0.98%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.86'
$ws.Range('E27').Value = '  -1.69%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.06'
$ws.Range('E28').Value = '  -1.36%  '
$ws.Range('E29').Value = '  +1.52%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '30.43'
$ws.Range('E30').Value = '  +0.35%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.57'
$ws.Range('E31').Value = '  -1.74%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '607.48'
$ws.Range('E32').Value = '  -6.59%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.76'
$ws.Range('E33').Value = '  +0.24%  '
$ws.Range('E34').Value = '  +0.37%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '59.83'
$ws.Range('E35').Value = '  +0.81%  '
$ws.Range('E36').Value = '  +8.44%  '
$ws.Range('B37').Value = 'Stacks'
$ws.Range('C37').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.69'
$ws.Range('E37').Value = '  +10.81%  '
$ws.Range('B38').Value = 'PEPE'
$ws.Range('C38').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D38').Value = '0.0₃0803'
$ws.Range('E38').Value = '  -2.09%  '
$ws.Range('E39').Value = '  -0.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '37.30'
$ws.Range('E40').Value = '  -3.83%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.384'
$ws.Range('E41').Value = '  -1.58%  '
$ws.Range('D42').Value = '3.264.03'
$ws.Range('E42').Value = '  +7.85%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  +0.30%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.00'
$ws.Range('E44').Value = '  +4.45%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.41'
$ws.Range('E45').Value = '  +1.20%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.55'
$ws.Range('E46').Value = '  -3.51%  '
$ws.Range('E47').Value = '  +0.76%  '
$ws.Range('E48').Value = '  -5.20%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.132'
$ws.Range('E49').Value = '  +1.28%  '
$ws.Range('E50').Value = '  -0.49%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '139.99'
$ws.Range('E51').Value = '  -0.65%  '
